# corrección de error behavior en los test case
# Update the Hogar claim test-case data row (row 2): new policy number and new
# claim date, then leave the selection where the test tool left it (E5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NroPoliza (E2): 11111003260 -> 11111003346 (text, keeps trailing space like row 3)
$ws.Range("E2").Value = "'11111003346 "

# FechaSiniestro (G2): 20/05/2021 -> 22/06/2021 (kept as text, quote-prefixed like before)
$ws.Range("G2").Value = "'22/06/2021"

# Cursor/selection ends up on E5 after the edits
$ws.Range("E5").Select() | Out-Null
